$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.881.25"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "2.306.63"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'307.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("D6").Value = "'96.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E7").Value = "  -2.05%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").Value = "'35.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "'18.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.64%  "
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "2.664.42"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "2.305.25"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "42.802.85"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "'13.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.04%  "
$ws.Range("D20").Value = "0.0₃0900"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "'6.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("D22").Value = "'67.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").Value = "'236.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").Value = "'2.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'4.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "'25.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D30").Value = "'166.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").Value = "'9.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'33.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'4.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").Value = "'5.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("D36").Value = "'17.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").Value = "'0.0695"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("D43").Value = "2.011.74"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("D45").Value = "'18.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.37%  "
$ws.Range("D46").Value = "'10.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("E47").Value = "  -6.75%  "
$ws.Range("D48").Value = "'2.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").Value = "'2.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.82%  "
$ws.Range("D50").Value = "'53.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").Value = "2.529.14"
$ws.Range("E51").Value = "  -0.01%  "
